$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for the "conv_v2" (apc_d3) strategy block (rows 117-149) ---
# New strategy columns J, K (added first so new shared strings land in the
# same table order the workbook ends up with)
$ws.Range("J118").Value = "W_Row_pipeline_ap_v2"
$ws.Range("K118").Value = "W_Col_pipeline_ap_v2"
# G118/H118 get renamed (WRp_ap -> W_Row_pipeline_ap_d3_c, WCp_ap -> W_Col_pipeline_ap_d3_c)
$ws.Range("G118").Value = "W_Row_pipeline_ap_d3_c"
$ws.Range("H118").Value = "W_Col_pipeline_ap_d3_c"
$ws.Range("I118").Value = "Filter_2_p_ap_v2"

# Directive labels for the new columns
$ws.Range("I121").Value = "pipeline"
$ws.Range("J122").Value = "pipeline"
$ws.Range("K123").Value = "pipeline"

$ws.Range("I125").Value = "apc_d3"
$ws.Range("J125").Value = "apc_d3"
$ws.Range("K125").Value = "apc_d3"

# Latency row
$ws.Range("I128").Value = 9745
$ws.Range("J128").Value = 172305
$ws.Range("K128").Value = 197473

# Resource usage rows (BRAM_18K, DSP48E, FF, LUT)
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 20
$ws.Range("K132").Value = 6

$ws.Range("I133").Value = 58
$ws.Range("J133").Value = 25
$ws.Range("K133").Value = 7

$ws.Range("I134").Value = 14582
$ws.Range("J134").Value = 3337
$ws.Range("K134").Value = 1421

$ws.Range("I135").Value = 13085
$ws.Range("J135").Value = 4661
$ws.Range("K135").Value = 1822

# SUM row
$ws.Range("I137").Formula = "=SUM(I132:I135)"
$ws.Range("J137").Formula = "=SUM(J132:J135)"
$ws.Range("K137").Formula = "=SUM(K132:K135)"

# Percent-of-total rows
$ws.Range("I140").Formula = "=I132/B132*100"
$ws.Range("J140").Formula = "=J132/B132*100"
$ws.Range("K140").Formula = "=K132/B132*100"

$ws.Range("I141").Formula = "=I133/B133*100"
$ws.Range("J141").Formula = "=J133/B133*100"
$ws.Range("K141").Formula = "=K133/B133*100"

$ws.Range("I142").Formula = "=I134/B134*100"
$ws.Range("J142").Formula = "=J134/B134*100"
$ws.Range("K142").Formula = "=K134/B134*100"

$ws.Range("I143").Formula = "=I135/B135*100"
$ws.Range("J143").Formula = "=J135/B135*100"
$ws.Range("K143").Formula = "=K135/B135*100"

# V row (average of the 4 percentages)
$ws.Range("I144").Formula = "=SUM(I140:I143)/4"
$ws.Range("J144").Formula = "=SUM(J140:J143)/4"
$ws.Range("K144").Formula = "=SUM(K140:K143)/4"
# J144/K144 pick up the same (10pt) font formatting used across the rest of row 144
$ws.Range("J144").Font.Size = $ws.Range("F144").Font.Size
$ws.Range("K144").Font.Size = $ws.Range("F144").Font.Size

# V2/V1 row
$ws.Range("I145").Formula = "=I144/C144"
$ws.Range("J145").Formula = "=J144/C144"
$ws.Range("K145").Formula = "=K144/C144"

# SpeedUp row
$ws.Range("I148").Formula = "=(C128/I128)"
$ws.Range("J148").Formula = "=(C128/J128)"
$ws.Range("K148").Formula = "=(C128/K128)"

# E row
$ws.Range("I149").Formula = "=I148/I145"
$ws.Range("J149").Formula = "=J148/J145"
$ws.Range("K149").Formula = "=K148/K145"

# --- Column widths to fit the longer new headers ---
# (values are pre-compensated for the host's automatic ~5/6-character
#  padding added on top of whatever is assigned to ColumnWidth)
$ws.Columns.Item(7).ColumnWidth = 25.166666666666668
$ws.Columns.Item(8).ColumnWidth = 24.0
$ws.Columns.Item(9).ColumnWidth = 25.166666666666668
$ws.Columns.Item(10).ColumnWidth = 23.333333333333336
$ws.Columns.Item(11).ColumnWidth = 22.333333333333336

# --- View state: selection moved to the newly filled area ---
$excel.ActiveWindow.ScrollRow = 125
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("I149").Select()
